$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ? pair): Target cluster changes from ECs to MuSCs, values refreshed with new TPM data
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 0.046374
$ws.Range("H2").Value = 0.139122
$ws.Range("I2").Value = 0.6592866045237633
$ws.Range("J2").Value = 0.6592866045237632
$ws.Range("M2").Value = 0.0002903333333333334
$ws.Range("N2").Value = 0.000871
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.000013463918
$ws.Range("R2").Value = 0.000121175262
$ws.Range("S2").Value = 0.6592866045237633
$ws.Range("T2").Value = 0.6592866045237632

# Row 3: Sending cluster changes from ECs to MuSCs (Target cluster stays MuSCs)
$ws.Range("A3").Value = "MuSCs"
$ws.Range("G3").Value = 0.02396566666666667
$ws.Range("H3").Value = 0.071897
$ws.Range("I3").Value = 0.3407133954762367
$ws.Range("J3").Value = 0.3407133954762367
$ws.Range("M3").Value = 0.0002903333333333334
$ws.Range("N3").Value = 0.000871
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.00000695803188888889
$ws.Range("R3").Value = 0.000062622287
$ws.Range("S3").Value = 0.3407133954762367
$ws.Range("T3").Value = 0.3407133954762367

$wb.Save()
